$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.627.37"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").Value = "1.845.81"
$ws.Range("E3").Value = "  +0.23%  "
$c = $ws.Range("D4")
$c.Value = "'1.034"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.28%  "
$c = $ws.Range("D5")
$c.Value = "'321.37"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.62%  "
$c = $ws.Range("D6")
$c.Value = "'1.029"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.19%  "
$c = $ws.Range("D7")
$c.Value = "'0.4375"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.03%  "
$c = $ws.Range("D8")
$c.Value = "'0.3784"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +1.15%  "
$c = $ws.Range("D9")
$c.Value = "'0.07367"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.23%  "
$c = $ws.Range("D10")
$c.Value = "'0.8797"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +0.46%  "
$c = $ws.Range("D11")
$c.Value = "'21.49"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("D12").Value = "1.859.10"
$ws.Range("E12").Value = "  +0.93%  "
$c = $ws.Range("D13")
$c.Value = "'5.484"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -0.13%  "
$c = $ws.Range("D14")
$c.Value = "'6.684"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +0.12%  "
$c = $ws.Range("D15")
$c.Value = "'0.07147"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -0.17%  "
$c = $ws.Range("D16")
$c.Value = "'84.79"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +2.48%  "
$ws.Range("E17").Value = "  +0.32%  "
$c = $ws.Range("D18")
$c.Value = "'0.000009033"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.41%  "
$c = $ws.Range("D19")
$c.Value = "'1.029"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.29%  "
$c = $ws.Range("D20")
$c.Value = "'15.41"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("D21").Value = "27.641.17"
$ws.Range("E21").Value = "  +0.37%  "
$c = $ws.Range("D22")
$c.Value = "'5.283"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.58%  "
$c = $ws.Range("D23")
$c.Value = "'11.27"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.95%  "
$ws.Range("D24").Value = "2.083.61"
$ws.Range("E24").Value = "  +1.12%  "
$c = $ws.Range("D25")
$c.Value = "'2.060"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +7.16%  "
$c = $ws.Range("D26")
$c.Value = "'158.28"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +0.58%  "
$c = $ws.Range("D27")
$c.Value = "'18.63"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.47%  "
$c = $ws.Range("D28")
$c.Value = "'1.984"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +2.51%  "
$c = $ws.Range("D29")
$c.Value = "'5.300"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +0.98%  "
$c = $ws.Range("D30")
$c.Value = "'117.45"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +1.14%  "
$c = $ws.Range("D31")
$c.Value = "'0.09029"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.67%  "
$c = $ws.Range("D32")
$c.Value = "'0.7681"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +0.35%  "
$c = $ws.Range("D33")
$c.Value = "'1.202"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.41%  "
$c = $ws.Range("D34")
$c.Value = "'2.998"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +4.15%  "
$c = $ws.Range("D35")
$c.Value = "'4.542"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.97%  "
$c = $ws.Range("D36")
$c.Value = "'1.031"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +0.04%  "
$c = $ws.Range("D37")
$c.Value = "'1.147"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +0.19%  "
$c = $ws.Range("D38")
$c.Value = "'0.01965"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.45%  "
$c = $ws.Range("D39")
$c.Value = "'0.05255"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +0.02%  "
$c = $ws.Range("D40")
$c.Value = "'2.838"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +1.75%  "
$c = $ws.Range("D41")
$c.Value = "'0.5157"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.14%  "
$c = $ws.Range("D42")
$c.Value = "'0.1664"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.18%  "
$c = $ws.Range("D43")
$c.Value = "'6.842"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +2.84%  "
$c = $ws.Range("D44")
$c.Value = "'8.692"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +2.11%  "
$c = $ws.Range("D45")
$c.Value = "'109.90"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +1.04%  "
$c = $ws.Range("D46")
$c.Value = "'10.69"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +1.03%  "
$c = $ws.Range("D47")
$c.Value = "'0.06600"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +3.98%  "
$c = $ws.Range("D48")
$c.Value = "'1.031"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +0.05%  "
$c = $ws.Range("D49")
$c.Value = "'1.693"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -0.62%  "
$c = $ws.Range("D50")
$c.Value = "'0.4681"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.79%  "
$c = $ws.Range("D51")
$c.Value = "'1.880"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.83%  "
